# Updated cryptos list (price + 1h volume change) per the scraper run.
# For D-column values that look like plain decimals (e.g. "305.83"), the
# cell is forced to text ("@") before assignment so Excel doesn't convert
# the string into a number, then the format is reset to "Normal" so no
# stray explicit style sticks on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.927.85'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '2.307.11'
$ws.Range("E3").Value = '  -0.15%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.33'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.98%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.512'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.63%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.504'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.84%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.70'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("E11").Value = '  +0.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '18.37'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.30%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.119'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.11%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("D15").Value = '2.662.93'
$ws.Range("E15").Value = '  -0.12%  '
$ws.Range("D16").Value = '2.300.16'
$ws.Range("E16").Value = '  -1.73%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.784'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.13%  '
$ws.Range("D18").Value = '42.839.74'
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.14'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = '0.0₃0904'
$ws.Range("E20").Value = '  -0.73%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.05'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.67'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.10%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.09'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.93%  '
$ws.Range("E26").Value = '  +0.98%  '
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.18%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '167.49'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.46%  '
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '33.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.80%  '
$ws.Range("E33").Value = '  +0.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.02'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.70%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.13%  '
$ws.Range("E37").Value = '  -0.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0695'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.77'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.57%  '
$ws.Range("E41").Value = '  -1.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.73'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("D43").Value = '2.007.61'
$ws.Range("E43").Value = '  -0.19%  '
$ws.Range("E44").Value = '  -2.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '18.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.63%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.02'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.97%  '
$ws.Range("E47").Value = '  -3.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.37%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.95'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +7.92%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '54.03'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.96%  '
$ws.Range("D51").Value = '2.530.65'
$ws.Range("E51").Value = '  -0.26%  '
